$d = $word.ActiveDocument

# Locate the paragraph that currently carries the "_GoBack" bookmark
# (this is the paragraph that should receive the new heading text).
$bk = $d.Bookmarks("_GoBack")
$targetRange = $d.Range($bk.Range.Start, $bk.Range.Start)
[void]$targetRange.Expand(4)  # wdParagraph

# Replace that paragraph's contents: drop the bookmark + rFonts hint,
# keep the <w:cs/> run-property, and add the new text run.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:rPr><w:cs/></w:rPr></w:pPr>' +
       '<w:r><w:t>How to upload photo process</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'
[void]$targetRange.InsertXML($xml)

# Move the "_GoBack" bookmark down to the last (empty) paragraph of the
# document body, right before the section properties.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$d.Bookmarks.Add("_GoBack", $lastPara.Range)
